$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header label in B2 from "unnamed: 1_level_1" to "total"
$ws.Range("B2").Value = "total"

# Remove the two empty sub-header separator rows ("situação do domicílio" and
# "grandes regiões e unidades da federação"), which shifts the data below them
# up and compacts the table (dimension goes from A1:I40 to A1:I38).
$ws.Rows("5").Delete()
$ws.Rows("7").Delete()

Write-Host ("UsedRange: " + $ws.UsedRange.Address())
